$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column widths between column A and column B
# (target stored widths are 13.7109375 / 14.42578125; the COM ColumnWidth
# model in this runtime only supports 1/6-character granularity, so we
# pick the ColumnWidth inputs that land in the middle of the bucket
# nearest to each target, i.e. the closest values reachable: 13.666... and 14.5)
$ws.Columns.Item(1).ColumnWidth = 12.875
$ws.Columns.Item(2).ColumnWidth = 13.666666666666664

# Update cell values
$ws.Range("A1").Value = -0.0324949512957039
$ws.Range("B1").Value = 0.032494951036158966

$ws.Range("A2").Value = 0.039459575608928654
$ws.Range("B2").Value = -0.039459575888081022

$ws.Range("A3").Value = -0.045257629914363662
$ws.Range("B3").Value = 0.045257629637817445
